$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update data values and font colors (negative => red, non-negative => green)
# matching the swapped font-color definitions (fontId 2 -> red, fontId 3 -> green).
$ws.Cells.Item(2, 3).Value = -0.35
$ws.Cells.Item(2, 3).Font.Color = 255
$ws.Cells.Item(2, 4).Value = 0.98
$ws.Cells.Item(2, 4).Font.Color = 32768
$ws.Cells.Item(2, 5).Value = -0.49
$ws.Cells.Item(2, 5).Font.Color = 255
$ws.Cells.Item(2, 6).Value = 4.52
$ws.Cells.Item(2, 6).Font.Color = 32768
$ws.Cells.Item(2, 7).Value = 0.46
$ws.Cells.Item(2, 7).Font.Color = 32768
$ws.Cells.Item(2, 8).Value = -3.42
$ws.Cells.Item(2, 8).Font.Color = 255

$ws.Cells.Item(3, 3).Value = -0.92
$ws.Cells.Item(3, 3).Font.Color = 255
$ws.Cells.Item(3, 4).Value = 0.17
$ws.Cells.Item(3, 4).Font.Color = 32768
$ws.Cells.Item(3, 5).Value = 3.94
$ws.Cells.Item(3, 5).Font.Color = 32768
$ws.Cells.Item(3, 6).Value = 15.05
$ws.Cells.Item(3, 6).Font.Color = 32768
$ws.Cells.Item(3, 7).Value = 12.96
$ws.Cells.Item(3, 7).Font.Color = 32768
$ws.Cells.Item(3, 8).Value = 3.42
$ws.Cells.Item(3, 8).Font.Color = 32768

$ws.Cells.Item(4, 3).Value = -0.34
$ws.Cells.Item(4, 3).Font.Color = 255
$ws.Cells.Item(4, 4).Value = 0.88
$ws.Cells.Item(4, 4).Font.Color = 32768
$ws.Cells.Item(4, 5).Value = 0.06
$ws.Cells.Item(4, 5).Font.Color = 32768
$ws.Cells.Item(4, 6).Value = 5.58
$ws.Cells.Item(4, 6).Font.Color = 32768
$ws.Cells.Item(4, 7).Value = 1.67
$ws.Cells.Item(4, 7).Font.Color = 32768
$ws.Cells.Item(4, 8).Value = -3.07
$ws.Cells.Item(4, 8).Font.Color = 255

$ws.Cells.Item(5, 3).Value = -0.65
$ws.Cells.Item(5, 3).Font.Color = 255
$ws.Cells.Item(5, 4).Value = -0.21
$ws.Cells.Item(5, 4).Font.Color = 255
$ws.Cells.Item(5, 5).Value = 3.32
$ws.Cells.Item(5, 5).Font.Color = 32768
$ws.Cells.Item(5, 6).Value = 14.15
$ws.Cells.Item(5, 6).Font.Color = 32768
$ws.Cells.Item(5, 7).Value = 14.2
$ws.Cells.Item(5, 7).Font.Color = 32768
$ws.Cells.Item(5, 8).Value = 3.55
$ws.Cells.Item(5, 8).Font.Color = 32768

$ws.Cells.Item(6, 3).Value = -0.39
$ws.Cells.Item(6, 3).Font.Color = 255
$ws.Cells.Item(6, 4).Value = 0.05
$ws.Cells.Item(6, 4).Font.Color = 32768
$ws.Cells.Item(6, 5).Value = 3.36
$ws.Cells.Item(6, 5).Font.Color = 32768
$ws.Cells.Item(6, 6).Value = 12.21
$ws.Cells.Item(6, 6).Font.Color = 32768
$ws.Cells.Item(6, 7).Value = 15.16
$ws.Cells.Item(6, 7).Font.Color = 32768
$ws.Cells.Item(6, 8).Value = 4.73
$ws.Cells.Item(6, 8).Font.Color = 32768

$ws.Cells.Item(7, 3).Value = -0.52
$ws.Cells.Item(7, 3).Font.Color = 255
$ws.Cells.Item(7, 4).Value = 0.97
$ws.Cells.Item(7, 4).Font.Color = 32768
$ws.Cells.Item(7, 5).Value = 1.28
$ws.Cells.Item(7, 5).Font.Color = 32768
$ws.Cells.Item(7, 6).Value = 9.43
$ws.Cells.Item(7, 6).Font.Color = 32768
$ws.Cells.Item(7, 7).Value = 2.45
$ws.Cells.Item(7, 7).Font.Color = 32768
$ws.Cells.Item(7, 8).Value = -1.82
$ws.Cells.Item(7, 8).Font.Color = 255

$ws.Cells.Item(8, 3).Value = -0.94
$ws.Cells.Item(8, 3).Font.Color = 255
$ws.Cells.Item(8, 4).Value = 0.55
$ws.Cells.Item(8, 4).Font.Color = 32768
$ws.Cells.Item(8, 5).Value = 2.64
$ws.Cells.Item(8, 5).Font.Color = 32768
$ws.Cells.Item(8, 6).Value = 12.46
$ws.Cells.Item(8, 6).Font.Color = 32768
$ws.Cells.Item(8, 7).Value = 7.35
$ws.Cells.Item(8, 7).Font.Color = 32768
$ws.Cells.Item(8, 8).Value = 0.02
$ws.Cells.Item(8, 8).Font.Color = 32768

$ws.Cells.Item(9, 3).Value = -0.2
$ws.Cells.Item(9, 3).Font.Color = 255
$ws.Cells.Item(9, 4).Value = 0.45
$ws.Cells.Item(9, 4).Font.Color = 32768
$ws.Cells.Item(9, 5).Value = 2.08
$ws.Cells.Item(9, 5).Font.Color = 32768
$ws.Cells.Item(9, 6).Value = 8.21
$ws.Cells.Item(9, 6).Font.Color = 32768
$ws.Cells.Item(9, 7).Value = 9.99
$ws.Cells.Item(9, 7).Font.Color = 32768
$ws.Cells.Item(9, 8).Value = 2.67
$ws.Cells.Item(9, 8).Font.Color = 32768

$ws.Cells.Item(10, 3).Value = -0.64
$ws.Cells.Item(10, 3).Font.Color = 255
$ws.Cells.Item(10, 4).Value = 0.45
$ws.Cells.Item(10, 4).Font.Color = 32768
$ws.Cells.Item(10, 5).Value = 3.21
$ws.Cells.Item(10, 5).Font.Color = 32768
$ws.Cells.Item(10, 6).Value = 20.03
$ws.Cells.Item(10, 6).Font.Color = 32768
$ws.Cells.Item(10, 7).Value = 21.86
$ws.Cells.Item(10, 7).Font.Color = 32768
$ws.Cells.Item(10, 8).Value = 3.98
$ws.Cells.Item(10, 8).Font.Color = 32768

$ws.Cells.Item(11, 3).Value = -0.08
$ws.Cells.Item(11, 3).Font.Color = 255
$ws.Cells.Item(11, 4).Value = 1.12
$ws.Cells.Item(11, 4).Font.Color = 32768
$ws.Cells.Item(11, 5).Value = 1.51
$ws.Cells.Item(11, 5).Font.Color = 32768
$ws.Cells.Item(11, 6).Value = -0.79
$ws.Cells.Item(11, 6).Font.Color = 255
$ws.Cells.Item(11, 7).Value = -5.83
$ws.Cells.Item(11, 7).Font.Color = 255
$ws.Cells.Item(11, 8).Value = -2.42
$ws.Cells.Item(11, 8).Font.Color = 255

$ws.Cells.Item(12, 3).Value = -1.11
$ws.Cells.Item(12, 3).Font.Color = 255
$ws.Cells.Item(12, 4).Value = 0.03
$ws.Cells.Item(12, 4).Font.Color = 32768
$ws.Cells.Item(12, 5).Value = 3.75
$ws.Cells.Item(12, 5).Font.Color = 32768
$ws.Cells.Item(12, 6).Value = 16.34
$ws.Cells.Item(12, 6).Font.Color = 32768
$ws.Cells.Item(12, 7).Value = 14.68
$ws.Cells.Item(12, 7).Font.Color = 32768
$ws.Cells.Item(12, 8).Value = 1.81
$ws.Cells.Item(12, 8).Font.Color = 32768

$ws.Cells.Item(13, 3).Value = 0.05
$ws.Cells.Item(13, 3).Font.Color = 32768
$ws.Cells.Item(13, 4).Value = 1.19
$ws.Cells.Item(13, 4).Font.Color = 32768
$ws.Cells.Item(13, 5).Value = 1.27
$ws.Cells.Item(13, 5).Font.Color = 32768
$ws.Cells.Item(13, 6).Value = 4.65
$ws.Cells.Item(13, 6).Font.Color = 32768
$ws.Cells.Item(13, 7).Value = 3.09
$ws.Cells.Item(13, 7).Font.Color = 32768
$ws.Cells.Item(13, 8).Value = 0.95
$ws.Cells.Item(13, 8).Font.Color = 32768

